$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 59
$ws.Range("A59").Value = "Federico Speroni"
$ws.Range("B58").Copy()
$ws.Range("B59").PasteSpecial(-4122)
$ws.Range("B59").Value = (Get-Date -Year 2017 -Month 5 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C59").Value = 2
$ws.Range("D59").Value = "Sprint 3 - Integración BackEnd y FrontEnd"
$ws.Range("E59").Value = "Lectura de lo realizado por Bruno"

# Row 60
$ws.Range("A60").Value = "Federico Speroni"
$ws.Range("B58").Copy()
$ws.Range("B60").PasteSpecial(-4122)
$ws.Range("B60").Value = (Get-Date -Year 2017 -Month 5 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C60").Value = 8
$ws.Range("D60").Value = "Sprint 3 - Integración BackEnd y FrontEnd"
$ws.Range("E60").Value = "Alta Cliente e ingreso de cliente "

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 41
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C60").Select()

$wb.Save()
